$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values after repull/push of data
$updates = @{
    2  = 0
    4  = -5
    5  = -3
    7  = 5
    8  = 4
    9  = 6
    10 = 5
    12 = 9
    13 = 2
    15 = 7
    16 = -1
    17 = 2
    18 = 1
    19 = 3
    20 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
